# Updates cryptos list values (price/volume) per the latest scrape.
# Leading apostrophe forces Excel to store the value as text, matching the
# original inlineStr cell type (prevents "572.11" etc. turning into a number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.608.45"
$ws.Range("E2").Value = "'  -0.83%  "

$ws.Range("D3").Value = "'2.458.40"
$ws.Range("E3").Value = "'  -0.20%  "

$ws.Range("E4").Value = "'  +0.06%  "

$ws.Range("D5").Value = "'572.11"
$ws.Range("E5").Value = "'  -1.38%  "

$ws.Range("D6").Value = "'146.93"
$ws.Range("E6").Value = "'  +0.20%  "

$ws.Range("E7").Value = "'  +0.06%  "

$ws.Range("D8").Value = "'0.532"
$ws.Range("E8").Value = "'  -1.79%  "

$ws.Range("D9").Value = "'0.111"
$ws.Range("E9").Value = "'  -0.29%  "

$ws.Range("E10").Value = "'  -0.33%  "

$ws.Range("D11").Value = "'5.28"
$ws.Range("E11").Value = "'  -0.28%  "

$ws.Range("D12").Value = "'0.351"
$ws.Range("E12").Value = "'  -1.19%  "

$ws.Range("D13").Value = "'28.75"
$ws.Range("E13").Value = "'  -1.50%  "

$ws.Range("D14").Value = "'0.0000175"
$ws.Range("E14").Value = "'  -2.55%  "

$ws.Range("D15").Value = "'2.905.54"
$ws.Range("E15").Value = "'  -0.15%  "

$ws.Range("D16").Value = "'62.626.41"
$ws.Range("E16").Value = "'  -0.60%  "

$ws.Range("D17").Value = "'2.458.65"
$ws.Range("E17").Value = "'  -0.29%  "

$ws.Range("E18").Value = "'  -0.50%  "

$ws.Range("D19").Value = "'10.86"
$ws.Range("E19").Value = "'  -2.43%  "

$ws.Range("D20").Value = "'324.71"
$ws.Range("E20").Value = "'  -1.82%  "

$ws.Range("E21").Value = "'  -0.13%  "

$ws.Range("E22").Value = "'  -3.57%  "

$ws.Range("E23").Value = "'  -0.06%  "

$ws.Range("D24").Value = "'9.97"
$ws.Range("E24").Value = "'  +10.62%  "

$ws.Range("D25").Value = "'65.32"
$ws.Range("E25").Value = "'  -1.82%  "

$ws.Range("D26").Value = "'640.10"
$ws.Range("E26").Value = "'  -4.20%  "

$ws.Range("D27").Value = "'2.586.31"
$ws.Range("E27").Value = "'  +0.23%  "

$ws.Range("D28").Value = "'0.0₃0967"
$ws.Range("E28").Value = "'  -4.13%  "

$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "'  -9.62%  "

$ws.Range("E30").Value = "'  -1.83%  "

$ws.Range("D31").Value = "'7.91"
$ws.Range("E31").Value = "'  -3.42%  "

$ws.Range("D32").Value = "'1.82"
$ws.Range("E32").Value = "'  -3.00%  "

$ws.Range("D33").Value = "'0.132"
$ws.Range("E33").Value = "'  -4.59%  "

$ws.Range("E34").Value = "'  -0.07%  "

$ws.Range("D35").Value = "'1.51"
$ws.Range("E35").Value = "'  -2.51%  "

$ws.Range("D36").Value = "'4.73"
$ws.Range("E36").Value = "'  -1.58%  "

$ws.Range("D37").Value = "'151.66"
$ws.Range("E37").Value = "'  -0.60%  "

$ws.Range("D38").Value = "'0.367"
$ws.Range("E38").Value = "'  -1.90%  "

$ws.Range("D39").Value = "'18.54"
$ws.Range("E39").Value = "'  -1.88%  "

$ws.Range("D40").Value = "'5.28"
$ws.Range("E40").Value = "'  -4.91%  "

$ws.Range("D41").Value = "'2.70"
$ws.Range("E41").Value = "'  -1.57%  "

$ws.Range("E42").Value = "'  -2.64%  "

$ws.Range("D43").Value = "'0.0₆0310"
$ws.Range("E43").Value = "'  -6.87%  "

$ws.Range("D45").Value = "'152.42"
$ws.Range("E45").Value = "'  +3.81%  "

$ws.Range("E46").Value = "'  +1.12%  "

$ws.Range("D47").Value = "'3.57"
$ws.Range("E47").Value = "'  -2.10%  "

$ws.Range("B48").Value = "'Mantle"
$ws.Range("C48").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.605"
$ws.Range("E48").Value = "'  -0.51%  "

$ws.Range("B49").Value = "'InjectiveProtocol"
$ws.Range("C49").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'20.20"
$ws.Range("E49").Value = "'  -2.89%  "

$ws.Range("D50").Value = "'0.0506"
$ws.Range("E50").Value = "'  -2.23%  "

$ws.Range("D51").Value = "'0.0908"
$ws.Range("E51").Value = "'  -1.65%  "
